$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Make "survey" the active sheet/tab (mirrors workbookView activeTab moving
# off "settings" and tabSelected shifting from the settings sheet to survey).
$ws.Activate() | Out-Null

# Insert a new row above the existing "send_sms" prompt row so a new
# "subject_name" text prompt is collected first.
$ws.Rows.Item(2).Insert() | Out-Null

# Fill in the new prompt row: type / name / display.text
# (set B2 before A2 so the new shared strings are appended in the same
# order as the target: subject_name, text, Enter the subject's name.)
$ws.Range("B2").Value = "subject_name"
$ws.Range("A2").Value = "text"
$ws.Range("C2").Value = "Enter the subject's name."

# Match the style used by the other "type"/header-ish cells on this sheet.
$ws.Range("A2:C2").Style = "Normal"
$ws.Rows.Item(2).RowHeight = 12

# Update the selection to land on the newly shifted send_sms row's value cell.
$ws.Range("D3").Select() | Out-Null
